# Use case "Escolher configuração ótima" - fill in the Pós condição (post-condition) cell
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C5").Value = "Configuração do carro completa"

# Reflect the cell that was last selected/active when the workbook was saved
$ws.Range("D9").Select()
